$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 134-135; Excel shifts existing rows 134:150 down to 136:152,
# carrying the D-column date style (s="2") down with them.
$ws.Rows("134:135").Insert()

# --- New row 134 (Primera / Caramelo, same market/region/product as its neighbours) ---
$ws.Cells.Item(134, 1).Value = 7
$ws.Cells.Item(134, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value = "Ñuble"
$ws.Cells.Item(134, 4).Value = 44474
$ws.Cells.Item(134, 5).Value = 16
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100108
$ws.Cells.Item(134, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(134, 9).Value = 100108005
$ws.Cells.Item(134, 10).Value = "Piña"
$ws.Cells.Item(134, 11).Value = "Caramelo"
$ws.Cells.Item(134, 12).Value = "Primera"
$ws.Cells.Item(134, 13).Value = 60
$ws.Cells.Item(134, 14).Value = 18000
$ws.Cells.Item(134, 15).Value = 19000
$ws.Cells.Item(134, 16).Value = 18500
$ws.Cells.Item(134, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(134, 18).Value = "Ecuador"
$ws.Cells.Item(134, 19).Value = 1542
$ws.Cells.Item(134, 20).Value = 12

# --- New row 135 (Segunda / Caramelo) ---
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 44474
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = "Fruta"
$ws.Cells.Item(135, 7).Value = 100108
$ws.Cells.Item(135, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(135, 9).Value = 100108005
$ws.Cells.Item(135, 10).Value = "Piña"
$ws.Cells.Item(135, 11).Value = "Caramelo"
$ws.Cells.Item(135, 12).Value = "Segunda"
$ws.Cells.Item(135, 13).Value = 60
$ws.Cells.Item(135, 14).Value = 18000
$ws.Cells.Item(135, 15).Value = 19000
$ws.Cells.Item(135, 16).Value = 18500
$ws.Cells.Item(135, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(135, 18).Value = "Ecuador"
$ws.Cells.Item(135, 19).Value = 1321
$ws.Cells.Item(135, 20).Value = 14
